$wb = $excel.ActiveWorkbook

# --- CENTRAL sheet: convert the "#1".."#76" row labels in column A from
# text to plain sequential numbers, and append a new summary row noting
# how many CENTRAL records were trials rather than reviews/protocols.
$central = $wb.Worksheets.Item("CENTRAL")

for ($i = 1; $i -le 76; $i++) {
    $central.Cells.Item($i, 1).Value = $i
}

# Copy the formatting of the existing "Total" row (76) down onto the new
# row 77 so the appended row matches the sheet's styling (A:C only, as D76
# already carries a footnote string that is not renumbered below it).
$central.Range("A76:C76").Copy()
$central.Range("A77:C77").PasteSpecial(-4122)  # xlPasteFormats

$central.Cells.Item(77, 2).Value = "Limiting to trials rather than reviews/protcols"
$central.Cells.Item(77, 3).Value = 1473

# --- View-state bookkeeping to mirror the authored workbook: CENTRAL
# becomes the active/selected sheet, scrolled down near the new row, with
# the cursor left on the newly-typed count.
$central.Activate()
$excel.ActiveWindow.ScrollRow = 61
$central.Range("B80").Select()

# WoS had also been scrolled down in the authored session.
$wos = $wb.Worksheets.Item("WoS")
$wos.Activate()
$excel.ActiveWindow.ScrollRow = 70

# Leave CENTRAL as the focused/active tab, matching the saved workbook.
$central.Activate()
